# Add a placeholder "confidence" column (N) to the event_type sheet with a
# default value of "unknown" for every data row, per the commit:
#   "Added placeholder confience column with default values of 'unknown'"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("event_type")

# Populate the data rows first so the "unknown" string is interned into the
# shared-strings table before "confidence" (matches the author's save order).
$ws.Range("N2:N74").Value = "unknown"

# Header cell for the new column.
$ws.Range("N1").Value = "confidence"

# Match the header's look-and-feel to the existing header row (style only,
# so it doesn't disturb the text we just wrote into N1).
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)

# Reflect the edit in the sheet's active selection, same as the source file.
$ws.Activate()
$ws.Range("N2:N74").Select()
